# Update the multiplication-problem table: each populated row (1, 5, 10,
# 15, 20) has 5 cells of "NNN×N=" text that get replaced with new values.
# We assign directly to each cell's Range.Text (rather than a global
# Find/Replace) because some new values coincide with other old values
# (e.g. "898×8=" -> "390×2=" while the original "390×2=" -> "162×5="),
# which would make a simple ordered find/replace ambiguous.

$d = $word.ActiveDocument
$table = $d.Tables.Item(1)

function Set-CellText($tbl, $row, $col, $text) {
    $tbl.Cell($row, $col).Range.Text = $text
}

# Row 1
Set-CellText $table 1 1 "459×2="
Set-CellText $table 1 2 "512×5="
Set-CellText $table 1 3 "151×3="
Set-CellText $table 1 4 "942×3="
Set-CellText $table 1 5 "921×5="

# Row 5
Set-CellText $table 5 1 "218×6="
Set-CellText $table 5 2 "968×4="
Set-CellText $table 5 3 "385×7="
Set-CellText $table 5 4 "598×6="
Set-CellText $table 5 5 "307×6="

# Row 10
Set-CellText $table 10 1 "293×3="
Set-CellText $table 10 2 "245×7="
Set-CellText $table 10 3 "162×5="
Set-CellText $table 10 4 "866×6="
Set-CellText $table 10 5 "255×8="

# Row 15
Set-CellText $table 15 1 "390×2="
Set-CellText $table 15 2 "294×9="
Set-CellText $table 15 3 "424×6="
Set-CellText $table 15 4 "761×6="
Set-CellText $table 15 5 "946×3="

# Row 20
Set-CellText $table 20 1 "436×4="
Set-CellText $table 20 2 "179×4="
Set-CellText $table 20 3 "699×6="
Set-CellText $table 20 4 "742×2="
Set-CellText $table 20 5 "890×7="

Write-Host "Done updating table cells"
